$d = $word.ActiveDocument

# "FECHA DE REALIZACIÓN DEL PROTOCOLO:" -> "FECHA DE REALIZACIÓN DEL CONSENTIMIENTO:"
# The word "PROTOCOLO" is replaced by "CONSENTIMIENTO" (same bold Book Antiqua 10pt run
# formatting is preserved automatically since Find/Replace keeps the original run's
# character formatting).
$d.Content.Find.Execute("PROTOCOLO", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CONSENTIMIENTO", 2)
